$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (O) to the table, mirroring the existing
# year-header (N4) and data (N5) cell formatting.
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value = 2021

$ws.Range("N3").Copy()
$ws.Range("O3").PasteSpecial(-4122)

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value = 4.0999999999999996

# Data corrections in the existing range.
$ws.Range("L5").Value = 1.6
$ws.Range("N5").Value = 3.1

# Move the active selection, matching the author's saved cursor position.
$ws.Range("P4").Select()
